# Registro de usuarios automatizado
$wb = $excel.ActiveWorkbook

# --- Remove Hoja2 and Hoja3 worksheets, keep only Hoja1 ---
$excel.DisplayAlerts = $false
foreach ($name in @("Hoja2", "Hoja3")) {
    $wb.Worksheets.Item($name).Delete() | Out-Null
}
$excel.DisplayAlerts = $true

# --- Work on Hoja1 ---
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# New data for column D (rows 1-20): number of credits / points per user
$dValues = @(58, 33, 994, 994, 1242, 973, 880, 1246, 375, 32, 501, 229, 1441, 975, 591, 387, 267, 47, 55, 246)
for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Column widths: C (3) and F (6)
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(6).ColumnWidth = 24.333333333333332

# View: scroll so row 3 is near the top, and select C12
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select() | Out-Null
